$wb = $excel.ActiveWorkbook

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Cells.Item(74, 8).Value = 38133.965  # H74: 53242.4 -> 38133.965
$ws.Cells.Item(74, 9).Value = 65421.438  # I74: 87174 -> 65421.438
$ws.Cells.Item(74, 10).Value = 1750.6666  # J74: 2345 -> 1750.6666
$ws.Cells.Item(74, 11).Value = 65421.438  # K74: 87174 -> 65421.438
$ws.Cells.Item(74, 12).Value = 1750.6666  # L74: 2345 -> 1750.6666
$ws.Cells.Item(74, 13).Value = -64547.438  # M74: -86300 -> -64547.438
$ws.Cells.Item(74, 14).Value = -3498.6666  # N74: -4093 -> -3498.6666
# Row 77
$ws.Cells.Item(77, 8).Value = 38133.965  # H77: 53242.4 -> 38133.965
$ws.Cells.Item(77, 9).Value = 65421.438  # I77: 87174 -> 65421.438
$ws.Cells.Item(77, 10).Value = 1750.6666  # J77: 2345 -> 1750.6666
$ws.Cells.Item(77, 11).Value = 327107.19  # K77: 435870 -> 327107.19
$ws.Cells.Item(77, 12).Value = 8753.333000000001  # L77: 11725 -> 8753.333000000001
$ws.Cells.Item(77, 13).Value = -322739.19  # M77: -431502 -> -322739.19
$ws.Cells.Item(77, 14).Value = -17489.333  # N77: -20461 -> -17489.333
# Row 132
$ws.Cells.Item(132, 8).Value = 1969675.9  # H132: 1727385.8 -> 1969675.9
$ws.Cells.Item(132, 9).Value = 2220721.8  # I132: 2375622 -> 2220721.8
$ws.Cells.Item(132, 10).Value = 919847.9399999999  # J132: 460378.8 -> 919847.9399999999
$ws.Cells.Item(132, 11).Value = 6662165.399999999  # K132: 7126866 -> 6662165.399999999
$ws.Cells.Item(132, 12).Value = 2759543.82  # L132: 1381136.4 -> 2759543.82
$ws.Cells.Item(132, 13).Value = -6659635.399999999  # M132: -7124336 -> -6659635.399999999
$ws.Cells.Item(132, 14).Value = -2764603.82  # N132: -1386196.4 -> -2764603.82

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 20727.475  # H134: 25014.447 -> 20727.475
$ws.Cells.Item(134, 9).Value = 1235  # I134: 1396.6487 -> 1235
$ws.Cells.Item(134, 10).Value = 93824.25  # J134: 112400.3 -> 93824.25
$ws.Cells.Item(134, 11).Value = 3705  # K134: 4189.9461 -> 3705
$ws.Cells.Item(134, 12).Value = 281472.75  # L134: 337200.9 -> 281472.75
$ws.Cells.Item(134, 13).Value = -1170  # M134: -1654.9461 -> -1170
$ws.Cells.Item(134, 14).Value = -286542.75  # N134: -342270.9 -> -286542.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Cells.Item(132, 8).Value = 1531.5555  # H132: 2203.9048 -> 1531.5555
$ws.Cells.Item(132, 9).Value = 972.4400000000001  # I132: 1497.6666 -> 972.4400000000001
$ws.Cells.Item(132, 10).Value = 2802.2727  # J132: 3145.5557 -> 2802.2727
$ws.Cells.Item(132, 11).Value = 2917.32  # K132: 4492.9998 -> 2917.32
$ws.Cells.Item(132, 12).Value = 8406.8181  # L132: 9436.667099999999 -> 8406.8181
$ws.Cells.Item(132, 13).Value = -387.3200000000002  # M132: -1962.9998 -> -387.3200000000002
$ws.Cells.Item(132, 14).Value = -13466.8181  # N132: -14496.6671 -> -13466.8181
# Row 134
$ws.Cells.Item(134, 8).Value = 1035.6897  # H134: 2057.3215 -> 1035.6897
$ws.Cells.Item(134, 9).Value = 899.0833  # I134: 1846.3125 -> 899.0833
$ws.Cells.Item(134, 10).Value = 1691.4  # J134: 2338.6667 -> 1691.4
$ws.Cells.Item(134, 11).Value = 2697.2499  # K134: 5538.9375 -> 2697.2499
$ws.Cells.Item(134, 12).Value = 5074.200000000001  # L134: 7016.000100000001 -> 5074.200000000001
$ws.Cells.Item(134, 13).Value = -162.2498999999998  # M134: -3003.9375 -> -162.2498999999998
$ws.Cells.Item(134, 14).Value = -10144.2  # N134: -12086.0001 -> -10144.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Cells.Item(62, 8).Value = 2050.6667  # H62: 2361.75 -> 2050.6667
$ws.Cells.Item(62, 9).Value = 1300  # I62: 0 -> 1300
$ws.Cells.Item(62, 10).Value = 2118.9092  # J62: 2361.75 -> 2118.9092
$ws.Cells.Item(62, 11).Value = 3900  # K62: 0 -> 3900
$ws.Cells.Item(62, 12).Value = 6356.7276  # L62: 7085.25 -> 6356.7276
$ws.Cells.Item(62, 13).Value = -3214  # M62: None -> -3214
$ws.Cells.Item(62, 14).Value = -7728.7276  # N62: -8457.25 -> -7728.7276
# Row 65
$ws.Cells.Item(65, 8).Value = 2050.6667  # H65: 2361.75 -> 2050.6667
$ws.Cells.Item(65, 9).Value = 1300  # I65: 0 -> 1300
$ws.Cells.Item(65, 10).Value = 2118.9092  # J65: 2361.75 -> 2118.9092
$ws.Cells.Item(65, 11).Value = 11700  # K65: 0 -> 11700
$ws.Cells.Item(65, 12).Value = 19070.1828  # L65: 21255.75 -> 19070.1828
$ws.Cells.Item(65, 13).Value = -8268  # M65: None -> -8268
$ws.Cells.Item(65, 14).Value = -25934.1828  # N65: -28119.75 -> -25934.1828
# Row 70
$ws.Cells.Item(70, 8).Value = 3400  # H70: 2977.889 -> 3400
$ws.Cells.Item(70, 9).Value = 3000  # I70: 1100.5 -> 3000
$ws.Cells.Item(70, 10).Value = 3480  # J70: 3514.2856 -> 3480
$ws.Cells.Item(70, 11).Value = 9000  # K70: 3301.5 -> 9000
$ws.Cells.Item(70, 12).Value = 10440  # L70: 10542.8568 -> 10440
$ws.Cells.Item(70, 13).Value = -8685  # M70: -2986.5 -> -8685
$ws.Cells.Item(70, 14).Value = -11070  # N70: -11172.8568 -> -11070
# Row 73
$ws.Cells.Item(73, 8).Value = 3400  # H73: 2977.889 -> 3400
$ws.Cells.Item(73, 9).Value = 3000  # I73: 1100.5 -> 3000
$ws.Cells.Item(73, 10).Value = 3480  # J73: 3514.2856 -> 3480
$ws.Cells.Item(73, 11).Value = 9000  # K73: 3301.5 -> 9000
$ws.Cells.Item(73, 12).Value = 10440  # L73: 10542.8568 -> 10440
$ws.Cells.Item(73, 13).Value = -7908  # M73: -2209.5 -> -7908
$ws.Cells.Item(73, 14).Value = -12624  # N73: -12726.8568 -> -12624
# Row 74
$ws.Cells.Item(74, 8).Value = 3739.8  # H74: 4163.8335 -> 3739.8
$ws.Cells.Item(74, 10).Value = 3739.8  # J74: 4163.8335 -> 3739.8
$ws.Cells.Item(74, 12).Value = 11219.4  # L74: 12491.5005 -> 11219.4
$ws.Cells.Item(74, 14).Value = -13341.4  # N74: -14613.5005 -> -13341.4
# Row 75
$ws.Cells.Item(75, 8).Value = 3557.5217  # H75: 5048.154 -> 3557.5217
$ws.Cells.Item(75, 9).Value = 2344.3333  # I75: 3266.5 -> 2344.3333
$ws.Cells.Item(75, 10).Value = 3739.5  # J75: 5372.091 -> 3739.5
$ws.Cells.Item(75, 11).Value = 7032.999899999999  # K75: 9799.5 -> 7032.999899999999
$ws.Cells.Item(75, 12).Value = 11218.5  # L75: 16116.273 -> 11218.5
$ws.Cells.Item(75, 13).Value = -6034.999899999999  # M75: -8801.5 -> -6034.999899999999
$ws.Cells.Item(75, 14).Value = -13214.5  # N75: -18112.273 -> -13214.5
# Row 76
$ws.Cells.Item(76, 8).Value = 14608.462  # H76: 13933 -> 14608.462
$ws.Cells.Item(76, 9).Value = 3500  # I76: 4233.3335 -> 3500
$ws.Cells.Item(76, 10).Value = 16628.182  # J76: 16578.363 -> 16628.182
$ws.Cells.Item(76, 11).Value = 10500  # K76: 12700.0005 -> 10500
$ws.Cells.Item(76, 12).Value = 49884.546  # L76: 49735.08900000001 -> 49884.546
$ws.Cells.Item(76, 13).Value = -10117  # M76: -12317.0005 -> -10117
$ws.Cells.Item(76, 14).Value = -50650.546  # N76: -50501.08900000001 -> -50650.546
# Row 77
$ws.Cells.Item(77, 8).Value = 3739.8  # H77: 4163.8335 -> 3739.8
$ws.Cells.Item(77, 10).Value = 3739.8  # J77: 4163.8335 -> 3739.8
$ws.Cells.Item(77, 12).Value = 33658.2  # L77: 37474.5015 -> 33658.2
$ws.Cells.Item(77, 14).Value = -44266.2  # N77: -48082.5015 -> -44266.2
# Row 78
$ws.Cells.Item(78, 8).Value = 3557.5217  # H78: 5048.154 -> 3557.5217
$ws.Cells.Item(78, 9).Value = 2344.3333  # I78: 3266.5 -> 2344.3333
$ws.Cells.Item(78, 10).Value = 3739.5  # J78: 5372.091 -> 3739.5
$ws.Cells.Item(78, 11).Value = 21098.9997  # K78: 29398.5 -> 21098.9997
$ws.Cells.Item(78, 12).Value = 33655.5  # L78: 48348.819 -> 33655.5
$ws.Cells.Item(78, 13).Value = -16106.9997  # M78: -24406.5 -> -16106.9997
$ws.Cells.Item(78, 14).Value = -43639.5  # N78: -58332.819 -> -43639.5
# Row 79
$ws.Cells.Item(79, 8).Value = 14608.462  # H79: 13933 -> 14608.462
$ws.Cells.Item(79, 9).Value = 3500  # I79: 4233.3335 -> 3500
$ws.Cells.Item(79, 10).Value = 16628.182  # J79: 16578.363 -> 16628.182
$ws.Cells.Item(79, 11).Value = 10500  # K79: 12700.0005 -> 10500
$ws.Cells.Item(79, 12).Value = 49884.546  # L79: 49735.08900000001 -> 49884.546
$ws.Cells.Item(79, 13).Value = -9174  # M79: -11374.0005 -> -9174
$ws.Cells.Item(79, 14).Value = -52536.546  # N79: -52387.08900000001 -> -52536.546
# Row 81
$ws.Cells.Item(81, 8).Value = 144287980  # H81: 125001230 -> 144287980
$ws.Cells.Item(81, 9).Value = 0  # I81: 368 -> 0
$ws.Cells.Item(81, 10).Value = 144287980  # J81: 333336000 -> 144287980
$ws.Cells.Item(81, 11).Value = 0  # K81: 1104 -> 0
$ws.Cells.Item(81, 12).Value = 432863940  # L81: 1000008000 -> 432863940
$ws.Cells.Item(81, 13).ClearContents()  # M81: clear (was 19)
$ws.Cells.Item(81, 14).Value = -432866186  # N81: -1000010246 -> -432866186
# Row 82
$ws.Cells.Item(82, 8).Value = 62502510  # H82: 62502260 -> 62502510
$ws.Cells.Item(82, 9).Value = 599  # I82: 787.75 -> 599
$ws.Cells.Item(82, 10).Value = 71431360  # J82: 83336080 -> 71431360
$ws.Cells.Item(82, 11).Value = 1797  # K82: 2363.25 -> 1797
$ws.Cells.Item(82, 12).Value = 214294080  # L82: 250008240 -> 214294080
$ws.Cells.Item(82, 13).Value = -1391  # M82: -1957.25 -> -1391
$ws.Cells.Item(82, 14).Value = -214294892  # N82: -250009052 -> -214294892
# Row 84
$ws.Cells.Item(84, 8).Value = 144287980  # H84: 125001230 -> 144287980
$ws.Cells.Item(84, 9).Value = 0  # I84: 368 -> 0
$ws.Cells.Item(84, 10).Value = 144287980  # J84: 333336000 -> 144287980
$ws.Cells.Item(84, 11).Value = 0  # K84: 3312 -> 0
$ws.Cells.Item(84, 12).Value = 1298591820  # L84: 3000024000 -> 1298591820
$ws.Cells.Item(84, 13).ClearContents()  # M84: clear (was 2304)
$ws.Cells.Item(84, 14).Value = -1298603052  # N84: -3000035232 -> -1298603052
# Row 85
$ws.Cells.Item(85, 8).Value = 62502510  # H85: 62502260 -> 62502510
$ws.Cells.Item(85, 9).Value = 599  # I85: 787.75 -> 599
$ws.Cells.Item(85, 10).Value = 71431360  # J85: 83336080 -> 71431360
$ws.Cells.Item(85, 11).Value = 1797  # K85: 2363.25 -> 1797
$ws.Cells.Item(85, 12).Value = 214294080  # L85: 250008240 -> 214294080
$ws.Cells.Item(85, 13).Value = -393  # M85: -959.25 -> -393
$ws.Cells.Item(85, 14).Value = -214296888  # N85: -250011048 -> -214296888
# Row 113
$ws.Cells.Item(113, 8).Value = 428.57574  # H113: 469.7143 -> 428.57574
$ws.Cells.Item(113, 9).Value = 398.92  # I113: 401.8 -> 398.92
$ws.Cells.Item(113, 10).Value = 521.25  # J113: 507.44446 -> 521.25
$ws.Cells.Item(113, 11).Value = 1196.76  # K113: 1205.4 -> 1196.76
$ws.Cells.Item(113, 12).Value = 1563.75  # L113: 1522.33338 -> 1563.75
$ws.Cells.Item(113, 13).Value = 973.24  # M113: 964.5999999999999 -> 973.24
$ws.Cells.Item(113, 14).Value = -5903.75  # N113: -5862.33338 -> -5903.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 20945.191  # H132: 42688.6 -> 20945.191
$ws.Cells.Item(132, 9).Value = 1479.6451  # I132: 2937.1428 -> 1479.6451
$ws.Cells.Item(132, 10).Value = 49680.047  # J132: 58147.5 -> 49680.047
$ws.Cells.Item(132, 11).Value = 4438.9353  # K132: 8811.428400000001 -> 4438.9353
$ws.Cells.Item(132, 12).Value = 149040.141  # L132: 174442.5 -> 149040.141
$ws.Cells.Item(132, 13).Value = -1908.9353  # M132: -6281.428400000001 -> -1908.9353
$ws.Cells.Item(132, 14).Value = -154100.141  # N132: -179502.5 -> -154100.141

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Cells.Item(132, 8).Value = 167247.67  # H132: 329193.72 -> 167247.67
$ws.Cells.Item(132, 9).Value = 43531.73  # I132: 147393.28 -> 43531.73
$ws.Cells.Item(132, 10).Value = 404782.28  # J132: 439854.88 -> 404782.28
$ws.Cells.Item(132, 11).Value = 130595.19  # K132: 442179.84 -> 130595.19
$ws.Cells.Item(132, 12).Value = 1214346.84  # L132: 1319564.64 -> 1214346.84
$ws.Cells.Item(132, 13).Value = -128065.19  # M132: -439649.84 -> -128065.19
$ws.Cells.Item(132, 14).Value = -1219406.84  # N132: -1324624.64 -> -1219406.84

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 4122.1943  # H132: 3588.0244 -> 4122.1943
$ws.Cells.Item(132, 9).Value = 977.44446  # I132: 963.2222 -> 977.44446
$ws.Cells.Item(132, 10).Value = 13556.444  # J132: 8650.143 -> 13556.444
$ws.Cells.Item(132, 11).Value = 2932.33338  # K132: 2889.6666 -> 2932.33338
$ws.Cells.Item(132, 12).Value = 40669.33199999999  # L132: 25950.429 -> 40669.33199999999
$ws.Cells.Item(132, 13).Value = -402.33338  # M132: -359.6666 -> -402.33338
$ws.Cells.Item(132, 14).Value = -45729.33199999999  # N132: -31010.429 -> -45729.33199999999
